$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 144.783305
$ws.Cells.Item(2, 8).Value = 434.349915
$ws.Cells.Item(2, 9).Value = 0.2430046335191003
$ws.Cells.Item(2, 10).Value = 0.251012682214973
$ws.Cells.Item(2, 13).Value = 0.6148836666666666
$ws.Cells.Item(2, 14).Value = 1.844651
$ws.Cells.Item(2, 15).Value = 0.1694384198480883
$ws.Cells.Item(2, 16).Value = 0.1898576456964817
$ws.Cells.Item(2, 17).Value = 89.02488945051832
$ws.Cells.Item(2, 18).Value = 801.2240050546649
$ws.Cells.Item(2, 19).Value = 0.04117432111924015
$ws.Cells.Item(2, 20).Value = 0.0476566768852939
$ws.Cells.Item(3, 7).Value = 144.783305
$ws.Cells.Item(3, 8).Value = 434.349915
$ws.Cells.Item(3, 9).Value = 0.2430046335191003
$ws.Cells.Item(3, 10).Value = 0.251012682214973
$ws.Cells.Item(3, 15).Value = 0.3591876696771482
$ws.Cells.Item(3, 16).Value = 0.4024738036936923
$ws.Cells.Item(3, 17).Value = 188.7213219626716
$ws.Cells.Item(3, 18).Value = 1698.491897664045
$ws.Cells.Item(3, 19).Value = 0.08728426803447506
$ws.Cells.Item(3, 20).Value = 0.1010260289864162
$ws.Cells.Item(4, 7).Value = 144.783305
$ws.Cells.Item(4, 8).Value = 434.349915
$ws.Cells.Item(4, 9).Value = 0.2430046335191003
$ws.Cells.Item(4, 10).Value = 0.251012682214973
$ws.Cells.Item(4, 13).Value = 0.1606106666666667
$ws.Cells.Item(4, 14).Value = 0.481832
$ws.Cells.Item(4, 15).Value = 0.04425815653597569
$ws.Cells.Item(4, 16).Value = 0.04959175971022552
$ws.Cells.Item(4, 17).Value = 23.25374313825333
$ws.Cells.Item(4, 18).Value = 209.28368824428
$ws.Cells.Item(4, 19).Value = 0.01075493710925575
$ws.Cells.Item(4, 20).Value = 0.01244816062062414
$ws.Cells.Item(5, 7).Value = 144.783305
$ws.Cells.Item(5, 8).Value = 434.349915
$ws.Cells.Item(5, 9).Value = 0.2430046335191003
$ws.Cells.Item(5, 10).Value = 0.251012682214973
$ws.Cells.Item(5, 13).Value = 1.170883
$ws.Cells.Item(5, 14).Value = 2.341766
$ws.Cells.Item(5, 15).Value = 0.322650569696364
$ws.Cells.Item(5, 16).Value = 0.2410223828421025
$ws.Cells.Item(5, 17).Value = 169.524310508315
$ws.Cells.Item(5, 18).Value = 1017.14586304989
$ws.Cells.Item(5, 19).Value = 0.07840558344379388
$ws.Cells.Item(5, 20).Value = 0.06049967479104023
$ws.Cells.Item(6, 7).Value = 144.783305
$ws.Cells.Item(6, 8).Value = 434.349915
$ws.Cells.Item(6, 9).Value = 0.2430046335191003
$ws.Cells.Item(6, 10).Value = 0.251012682214973
$ws.Cells.Item(6, 13).Value = 0.379099
$ws.Cells.Item(6, 14).Value = 1.137297
$ws.Cells.Item(6, 15).Value = 0.1044651842424238
$ws.Cells.Item(6, 16).Value = 0.1170544080574979
$ws.Cells.Item(6, 17).Value = 54.887206142195
$ws.Cells.Item(6, 18).Value = 493.9848552797549
$ws.Cells.Item(6, 19).Value = 0.02538552381233549
$ws.Cells.Item(6, 20).Value = 0.0293821409315985
$ws.Cells.Item(7, 7).Value = 82.24887099999999
$ws.Cells.Item(7, 9).Value = 0.1380466950572427
$ws.Cells.Item(7, 10).Value = 0.1425959278859072
$ws.Cells.Item(7, 13).Value = 0.6148836666666666
$ws.Cells.Item(7, 14).Value = 1.844651
$ws.Cells.Item(7, 15).Value = 0.1694384198480883
$ws.Cells.Item(7, 16).Value = 0.1898576456964817
$ws.Cells.Item(7, 17).Value = 50.57348737967366
$ws.Cells.Item(7, 18).Value = 455.1613864170629
$ws.Cells.Item(7, 19).Value = 0.0233904138757501
$ws.Cells.Item(7, 20).Value = 0.02707292715432363
$ws.Cells.Item(8, 7).Value = 82.24887099999999
$ws.Cells.Item(8, 9).Value = 0.1380466950572427
$ws.Cells.Item(8, 10).Value = 0.1425959278859072
$ws.Cells.Item(8, 15).Value = 0.3591876696771482
$ws.Cells.Item(8, 16).Value = 0.4024738036936923
$ws.Cells.Item(8, 18).Value = 964.8836306472989
$ws.Cells.Item(8, 19).Value = 0.04958467070424288
$ws.Cells.Item(8, 20).Value = 0.05739112548747251
$ws.Cells.Item(9, 7).Value = 82.24887099999999
$ws.Cells.Item(9, 9).Value = 0.1380466950572427
$ws.Cells.Item(9, 10).Value = 0.1425959278859072
$ws.Cells.Item(9, 13).Value = 0.1606106666666667
$ws.Cells.Item(9, 14).Value = 0.481832
$ws.Cells.Item(9, 15).Value = 0.04425815653597569
$ws.Cells.Item(9, 16).Value = 0.04959175971022552
$ws.Cells.Item(9, 17).Value = 13.21004600389067
$ws.Cells.Item(9, 18).Value = 118.890414035016
$ws.Cells.Item(9, 19).Value = 0.006109692239117549
$ws.Cells.Item(9, 20).Value = 0.007071582991374555
$ws.Cells.Item(10, 7).Value = 82.24887099999999
$ws.Cells.Item(10, 9).Value = 0.1380466950572427
$ws.Cells.Item(10, 10).Value = 0.1425959278859072
$ws.Cells.Item(10, 13).Value = 1.170883
$ws.Cells.Item(10, 14).Value = 2.341766
$ws.Cells.Item(10, 15).Value = 0.322650569696364
$ws.Cells.Item(10, 16).Value = 0.2410223828421025
$ws.Cells.Item(10, 17).Value = 96.30380482309299
$ws.Cells.Item(10, 18).Value = 577.8228289385579
$ws.Cells.Item(10, 19).Value = 0.04454084480491959
$ws.Cells.Item(10, 20).Value = 0.03436881032264196
$ws.Cells.Item(11, 7).Value = 82.24887099999999
$ws.Cells.Item(11, 9).Value = 0.1380466950572427
$ws.Cells.Item(11, 10).Value = 0.1425959278859072
$ws.Cells.Item(11, 13).Value = 0.379099
$ws.Cells.Item(11, 14).Value = 1.137297
$ws.Cells.Item(11, 15).Value = 0.1044651842424238
$ws.Cells.Item(11, 16).Value = 0.1170544080574979
$ws.Cells.Item(11, 17).Value = 31.180464747229
$ws.Cells.Item(11, 18).Value = 280.624182725061
$ws.Cells.Item(11, 19).Value = 0.01442107343321255
$ws.Cells.Item(11, 20).Value = 0.01669148193009453
$ws.Cells.Item(12, 7).Value = 163.8590903333333
$ws.Cells.Item(12, 8).Value = 491.577271
$ws.Cells.Item(12, 9).Value = 0.2750214756820535
$ws.Cells.Item(12, 10).Value = 0.284084617144743
$ws.Cells.Item(12, 13).Value = 0.6148836666666666
$ws.Cells.Item(12, 14).Value = 1.844651
$ws.Cells.Item(12, 15).Value = 0.1694384198480883
$ws.Cells.Item(12, 16).Value = 0.1898576456964817
$ws.Cells.Item(12, 17).Value = 100.7542782808246
$ws.Cells.Item(12, 18).Value = 906.7885045274209
$ws.Cells.Item(12, 19).Value = 0.04659920426385659
$ws.Cells.Item(12, 20).Value = 0.05393563658968727
$ws.Cells.Item(13, 7).Value = 163.8590903333333
$ws.Cells.Item(13, 8).Value = 491.577271
$ws.Cells.Item(13, 9).Value = 0.2750214756820535
$ws.Cells.Item(13, 10).Value = 0.284084617144743
$ws.Cells.Item(13, 15).Value = 0.3591876696771482
$ws.Cells.Item(13, 16).Value = 0.4024738036936923
$ws.Cells.Item(13, 17).Value = 213.5861185328481
$ws.Cells.Item(13, 18).Value = 1922.275066795633
$ws.Cells.Item(13, 19).Value = 0.09878432296140727
$ws.Cells.Item(13, 20).Value = 0.114336616433111
$ws.Cells.Item(14, 7).Value = 163.8590903333333
$ws.Cells.Item(14, 8).Value = 491.577271
$ws.Cells.Item(14, 9).Value = 0.2750214756820535
$ws.Cells.Item(14, 10).Value = 0.284084617144743
$ws.Cells.Item(14, 13).Value = 0.1606106666666667
$ws.Cells.Item(14, 14).Value = 0.481832
$ws.Cells.Item(14, 15).Value = 0.04425815653597569
$ws.Cells.Item(14, 16).Value = 0.04959175971022552
$ws.Cells.Item(14, 17).Value = 26.31751773783023
$ws.Cells.Item(14, 18).Value = 236.857659640472
$ws.Cells.Item(14, 19).Value = 0.01217194352149136
$ws.Cells.Item(14, 20).Value = 0.01408825607081351
$ws.Cells.Item(15, 7).Value = 163.8590903333333
$ws.Cells.Item(15, 8).Value = 491.577271
$ws.Cells.Item(15, 9).Value = 0.2750214756820535
$ws.Cells.Item(15, 10).Value = 0.284084617144743
$ws.Cells.Item(15, 13).Value = 1.170883
$ws.Cells.Item(15, 14).Value = 2.341766
$ws.Cells.Item(15, 15).Value = 0.322650569696364
$ws.Cells.Item(15, 16).Value = 0.2410223828421025
$ws.Cells.Item(15, 17).Value = 191.8598232667643
$ws.Cells.Item(15, 18).Value = 1151.158939600586
$ws.Cells.Item(15, 19).Value = 0.08873583580754928
$ws.Cells.Item(15, 20).Value = 0.06847075135301237
$ws.Cells.Item(16, 7).Value = 163.8590903333333
$ws.Cells.Item(16, 8).Value = 491.577271
$ws.Cells.Item(16, 9).Value = 0.2750214756820535
$ws.Cells.Item(16, 10).Value = 0.284084617144743
$ws.Cells.Item(16, 13).Value = 0.379099
$ws.Cells.Item(16, 14).Value = 1.137297
$ws.Cells.Item(16, 15).Value = 0.1044651842424238
$ws.Cells.Item(16, 16).Value = 0.1170544080574979
$ws.Cells.Item(16, 17).Value = 62.11881728627634
$ws.Cells.Item(16, 18).Value = 559.0693555764869
$ws.Cells.Item(16, 19).Value = 0.028730169127749
$ws.Cells.Item(16, 20).Value = 0.03325335669811882
$ws.Cells.Item(17, 7).Value = 57.0238095
$ws.Cells.Item(17, 8).Value = 114.047619
$ws.Cells.Item(17, 9).Value = 0.09570889357312636
$ws.Cells.Item(17, 10).Value = 0.06590860906562239
$ws.Cells.Item(17, 13).Value = 0.6148836666666666
$ws.Cells.Item(17, 14).Value = 1.844651
$ws.Cells.Item(17, 15).Value = 0.1694384198480883
$ws.Cells.Item(17, 16).Value = 0.1898576456964817
$ws.Cells.Item(17, 17).Value = 35.06300907266149
$ws.Cells.Item(17, 18).Value = 210.378054435969
$ws.Cells.Item(17, 19).Value = 0.01621676369243938
$ws.Cells.Item(17, 20).Value = 0.01251325334832886
$ws.Cells.Item(18, 7).Value = 57.0238095
$ws.Cells.Item(18, 8).Value = 114.047619
$ws.Cells.Item(18, 9).Value = 0.09570889357312636
$ws.Cells.Item(18, 10).Value = 0.06590860906562239
$ws.Cells.Item(18, 15).Value = 0.3591876696771482
$ws.Cells.Item(18, 16).Value = 0.4024738036936923
$ws.Cells.Item(18, 17).Value = 74.3290720721395
$ws.Cells.Item(18, 18).Value = 445.974432432837
$ws.Cells.Item(18, 19).Value = 0.03437745444990944
$ws.Cells.Item(18, 20).Value = 0.02652648858680162
$ws.Cells.Item(19, 7).Value = 57.0238095
$ws.Cells.Item(19, 8).Value = 114.047619
$ws.Cells.Item(19, 9).Value = 0.09570889357312636
$ws.Cells.Item(19, 10).Value = 0.06590860906562239
$ws.Cells.Item(19, 13).Value = 0.1606106666666667
$ws.Cells.Item(19, 14).Value = 0.481832
$ws.Cells.Item(19, 15).Value = 0.04425815653597569
$ws.Cells.Item(19, 16).Value = 0.04959175971022552
$ws.Cells.Item(19, 17).Value = 9.158632059668001
$ws.Cells.Item(19, 18).Value = 54.95179235800801
$ws.Cells.Item(19, 19).Value = 0.004235899193644465
$ws.Cells.Item(19, 20).Value = 0.003268523903617536
$ws.Cells.Item(20, 7).Value = 57.0238095
$ws.Cells.Item(20, 8).Value = 114.047619
$ws.Cells.Item(20, 9).Value = 0.09570889357312636
$ws.Cells.Item(20, 10).Value = 0.06590860906562239
$ws.Cells.Item(20, 13).Value = 1.170883
$ws.Cells.Item(20, 14).Value = 2.341766
$ws.Cells.Item(20, 15).Value = 0.322650569696364
$ws.Cells.Item(20, 16).Value = 0.2410223828421025
$ws.Cells.Item(20, 17).Value = 66.76820913878849
$ws.Cells.Item(20, 18).Value = 267.072836555154
$ws.Cells.Item(20, 19).Value = 0.03088052903637789
$ws.Cells.Item(20, 20).Value = 0.01588545000680491
$ws.Cells.Item(21, 7).Value = 57.0238095
$ws.Cells.Item(21, 8).Value = 114.047619
$ws.Cells.Item(21, 9).Value = 0.09570889357312636
$ws.Cells.Item(21, 10).Value = 0.06590860906562239
$ws.Cells.Item(21, 13).Value = 0.379099
$ws.Cells.Item(21, 14).Value = 1.137297
$ws.Cells.Item(21, 15).Value = 0.1044651842424238
$ws.Cells.Item(21, 16).Value = 0.1170544080574979
$ws.Cells.Item(21, 17).Value = 21.6176691576405
$ws.Cells.Item(21, 18).Value = 129.706014945843
$ws.Cells.Item(21, 19).Value = 0.009998247200755176
$ws.Cells.Item(21, 20).Value = 0.00771489322006947
$ws.Cells.Item(22, 7).Value = 147.8896333333333
$ws.Cells.Item(22, 8).Value = 443.6689
$ws.Cells.Item(22, 9).Value = 0.2482183021684772
$ws.Cells.Item(22, 10).Value = 0.2563981636887546
$ws.Cells.Item(22, 13).Value = 0.6148836666666666
$ws.Cells.Item(22, 14).Value = 1.844651
$ws.Cells.Item(22, 15).Value = 0.1694384198480883
$ws.Cells.Item(22, 16).Value = 0.1898576456964817
$ws.Cells.Item(22, 17).Value = 90.93492000598889
$ws.Cells.Item(22, 18).Value = 818.4142800538999
$ws.Cells.Item(22, 19).Value = 0.04205771689680209
$ws.Cells.Item(22, 20).Value = 0.04867915171884809
$ws.Cells.Item(23, 7).Value = 147.8896333333333
$ws.Cells.Item(23, 8).Value = 443.6689
$ws.Cells.Item(23, 9).Value = 0.2482183021684772
$ws.Cells.Item(23, 10).Value = 0.2563981636887546
$ws.Cells.Item(23, 15).Value = 0.3591876696771482
$ws.Cells.Item(23, 16).Value = 0.4024738036936923
$ws.Cells.Item(23, 17).Value = 192.7703412160778
$ws.Cells.Item(23, 18).Value = 1734.9330709447
$ws.Cells.Item(23, 19).Value = 0.08915695352711354
$ws.Cells.Item(23, 20).Value = 0.103193544199891
$ws.Cells.Item(24, 7).Value = 147.8896333333333
$ws.Cells.Item(24, 8).Value = 443.6689
$ws.Cells.Item(24, 9).Value = 0.2482183021684772
$ws.Cells.Item(24, 10).Value = 0.2563981636887546
$ws.Cells.Item(24, 13).Value = 0.1606106666666667
$ws.Cells.Item(24, 14).Value = 0.481832
$ws.Cells.Item(24, 15).Value = 0.04425815653597569
$ws.Cells.Item(24, 16).Value = 0.04959175971022552
$ws.Cells.Item(24, 17).Value = 23.75265260275556
$ws.Cells.Item(24, 18).Value = 213.7738734248
$ws.Cells.Item(24, 19).Value = 0.01098568447246658
$ws.Cells.Item(24, 20).Value = 0.01271523612379579
$ws.Cells.Item(25, 7).Value = 147.8896333333333
$ws.Cells.Item(25, 8).Value = 443.6689
$ws.Cells.Item(25, 9).Value = 0.2482183021684772
$ws.Cells.Item(25, 10).Value = 0.2563981636887546
$ws.Cells.Item(25, 13).Value = 1.170883
$ws.Cells.Item(25, 14).Value = 2.341766
$ws.Cells.Item(25, 15).Value = 0.322650569696364
$ws.Cells.Item(25, 16).Value = 0.2410223828421025
$ws.Cells.Item(25, 17).Value = 173.1614575462333
$ws.Cells.Item(25, 18).Value = 1038.9687452774
$ws.Cells.Item(25, 19).Value = 0.0800877766037234
$ws.Cells.Item(25, 20).Value = 0.06179769636860307
$ws.Cells.Item(26, 7).Value = 147.8896333333333
$ws.Cells.Item(26, 8).Value = 443.6689
$ws.Cells.Item(26, 9).Value = 0.2482183021684772
$ws.Cells.Item(26, 10).Value = 0.2563981636887546
$ws.Cells.Item(26, 13).Value = 0.379099
$ws.Cells.Item(26, 14).Value = 1.137297
$ws.Cells.Item(26, 15).Value = 0.1044651842424238
$ws.Cells.Item(26, 16).Value = 0.1170544080574979
$ws.Cells.Item(26, 17).Value = 56.06481210703334
$ws.Cells.Item(26, 18).Value = 504.5833089633
$ws.Cells.Item(26, 19).Value = 0.0259301706683716
$ws.Cells.Item(26, 20).Value = 0.03001253527761662
